$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Select column N (the "Late" column) and insert a new blank column before it,
# shifting "Late" and "Outstanding" one column to the right (N->O, O->P, P->Q).
$ws.Columns.Item(14).Select() | Out-Null
$ws.Columns.Item(14).Insert()

# Make "Repayment Schedule" the active sheet/tab (previously "Transactions" was active).
$ws.Activate()
